# Update "想去人数" (want-to-go count, column F) figures across sheets to
# match the freshly scraped numbers (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value  = 3969
$ws.Range("F7").Value  = 2614
$ws.Range("F13").Value = 484
$ws.Range("F19").Value = 495
$ws.Range("F20").Value = 701
$ws.Range("F26").Value = 201
$ws.Range("F30").Value = 4896
$ws.Range("F31").Value = 4523
$ws.Range("F34").Value = 108
$ws.Range("F39").Value = 19
$ws.Range("F40").Value = 543
$ws.Range("F42").Value = 1329

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 163
$ws.Range("F4").Value = 2371
$ws.Range("F5").Value = 133

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 163
$ws.Range("F8").Value  = 3969
$ws.Range("F9").Value  = 3969
$ws.Range("F10").Value = 2614
$ws.Range("F13").Value = 133
$ws.Range("F16").Value = 484
$ws.Range("F19").Value = 243
$ws.Range("F23").Value = 495
$ws.Range("F24").Value = 701
$ws.Range("F34").Value = 4896
$ws.Range("F35").Value = 4523
$ws.Range("F41").Value = 543
